# The deck ships with two themes:
#   ppt/theme/theme1.xml  -> used by the (only) Slide Master  -> currently "Integral"
#   ppt/theme/theme2.xml  -> used by the Notes Master          -> currently "Office Theme"
#
# The authored edit swaps the two themes' contents so that the Slide Master
# (and therefore every slide) now renders with the standard "Office Theme"
# palette, while the "Integral" palette moves over to where "Office Theme"
# used to live.
#
# Drive this through the Design/Theme color scheme on the presentation's
# slide master, setting each of the twelve theme colors (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) to the standard Office theme values.

function ToCOMColor([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Standard "Office Theme" color scheme values (the theme this deck's Notes
# Master already carries), applied here so the Slide Master/theme1 switches
# over to it.
$colorScheme.Colors(1).RGB  = ToCOMColor 0x00 0x00 0x00   # dk1
$colorScheme.Colors(2).RGB  = ToCOMColor 0xFF 0xFF 0xFF   # lt1
$colorScheme.Colors(3).RGB  = ToCOMColor 0x44 0x54 0x6A   # dk2
$colorScheme.Colors(4).RGB  = ToCOMColor 0xE7 0xE6 0xE6   # lt2
$colorScheme.Colors(5).RGB  = ToCOMColor 0x5B 0x9B 0xD5   # accent1
$colorScheme.Colors(6).RGB  = ToCOMColor 0xED 0x7D 0x31   # accent2
$colorScheme.Colors(7).RGB  = ToCOMColor 0xA5 0xA5 0xA5   # accent3
$colorScheme.Colors(8).RGB  = ToCOMColor 0xFF 0xC0 0x00   # accent4
$colorScheme.Colors(9).RGB  = ToCOMColor 0x44 0x72 0xC4   # accent5
$colorScheme.Colors(10).RGB = ToCOMColor 0x70 0xAD 0x47   # accent6
$colorScheme.Colors(11).RGB = ToCOMColor 0x05 0x63 0xC1   # hlink
$colorScheme.Colors(12).RGB = ToCOMColor 0x95 0x4F 0x72   # folHlink
